$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 10841
$ws.Cells.Item(2, 5).Value = 628
$ws.Cells.Item(2, 6).Value = 628
$ws.Cells.Item(2, 7).Value = 347
$ws.Cells.Item(2, 8).Value = 238
$ws.Cells.Item(2, 9).Value = 232
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(2, 11).Value = 10898
$ws.Cells.Item(2, 12).Value = 8251
$ws.Cells.Item(2, 13).Value = 2647
$ws.Cells.Item(2, 14).Value = 2221
$ws.Cells.Item(2, 15).Value = 426
$ws.Cells.Item(2, 16).Value = 78
$ws.Cells.Item(2, 17).Value = 883
$ws.Cells.Item(2, 18).Value = -328
$ws.Cells.Item(2, 19).Value = -477
$ws.Cells.Item(2, 20).Value = 373
$ws.Cells.Item(2, 21).Value = 510
$ws.Cells.Item(2, 22).Value = 4977
$ws.Cells.Item(2, 23).Value = 5.79
$ws.Cells.Item(2, 24).Value = 2.19
$ws.Cells.Item(2, 25).Value = 11
$ws.Cells.Item(2, 26).Value = 2.18
$ws.Cells.Item(2, 27).Value = 311.66
$ws.Cells.Item(2, 28).Value = 3882.43
$ws.Cells.Item(2, 29).Value = 1480
$ws.Cells.Item(2, 30).Value = 12.61
$ws.Cells.Item(2, 31).Value = 14670
$ws.Cells.Item(2, 32).Value = 1.27
$ws.Cells.Item(2, 33).Value = 100
$ws.Cells.Item(2, 34).Value = 0.54
$ws.Cells.Item(2, 35).Value = 6.56
$ws.Cells.Item(2, 36).Value = 14735240

# Row 3
$ws.Cells.Item(3, 4).Value = 12040
$ws.Cells.Item(3, 5).Value = 874
$ws.Cells.Item(3, 6).Value = 874
$ws.Cells.Item(3, 7).Value = 585
$ws.Cells.Item(3, 8).Value = 432
$ws.Cells.Item(3, 9).Value = 415
$ws.Cells.Item(3, 10).Value = 17
$ws.Cells.Item(3, 11).Value = 11216
$ws.Cells.Item(3, 12).Value = 8181
$ws.Cells.Item(3, 13).Value = 3034
$ws.Cells.Item(3, 14).Value = 2593
$ws.Cells.Item(3, 15).Value = 442
$ws.Cells.Item(3, 16).Value = 78
$ws.Cells.Item(3, 17).Value = 1126
$ws.Cells.Item(3, 18).Value = -564
$ws.Cells.Item(3, 19).Value = -496
$ws.Cells.Item(3, 20).Value = 557
$ws.Cells.Item(3, 21).Value = 569
$ws.Cells.Item(3, 22).Value = 4706
$ws.Cells.Item(3, 23).Value = 7.26
$ws.Cells.Item(3, 24).Value = 3.59
$ws.Cells.Item(3, 25).Value = 17.24
$ws.Cells.Item(3, 26).Value = 3.91
$ws.Cells.Item(3, 27).Value = 269.61
$ws.Cells.Item(3, 28).Value = 4352.84
$ws.Cells.Item(3, 29).Value = 2653
$ws.Cells.Item(3, 30).Value = 20.12
$ws.Cells.Item(3, 31).Value = 17121
$ws.Cells.Item(3, 32).Value = 3.12
$ws.Cells.Item(3, 33).Value = 200
$ws.Cells.Item(3, 34).Value = 0.37
$ws.Cells.Item(3, 35).Value = 7.31
$ws.Cells.Item(3, 36).Value = 14735240

# Row 4
$ws.Cells.Item(4, 4).Value = 11888
$ws.Cells.Item(4, 5).Value = 659
$ws.Cells.Item(4, 6).Value = 659
$ws.Cells.Item(4, 7).Value = 557
$ws.Cells.Item(4, 8).Value = 446
$ws.Cells.Item(4, 9).Value = 361
$ws.Cells.Item(4, 10).Value = 85
$ws.Cells.Item(4, 11).Value = 11778
$ws.Cells.Item(4, 12).Value = 7729
$ws.Cells.Item(4, 13).Value = 4049
$ws.Cells.Item(4, 14).Value = 2964
$ws.Cells.Item(4, 15).Value = 1084
$ws.Cells.Item(4, 16).Value = 78
$ws.Cells.Item(4, 17).Value = 821
$ws.Cells.Item(4, 18).Value = -665
$ws.Cells.Item(4, 19).Value = -146
$ws.Cells.Item(4, 20).Value = 677
$ws.Cells.Item(4, 21).Value = 144
$ws.Cells.Item(4, 22).Value = 4152
$ws.Cells.Item(4, 23).Value = 5.54
$ws.Cells.Item(4, 24).Value = 3.75
$ws.Cells.Item(4, 25).Value = 13
$ws.Cells.Item(4, 26).Value = 3.88
$ws.Cells.Item(4, 27).Value = 190.91
$ws.Cells.Item(4, 28).Value = 4752.53
$ws.Cells.Item(4, 29).Value = 2309
$ws.Cells.Item(4, 30).Value = 11.99
$ws.Cells.Item(4, 31).Value = 19512
$ws.Cells.Item(4, 32).Value = 1.42
$ws.Cells.Item(4, 33).Value = 200
$ws.Cells.Item(4, 34).Value = 0.72
$ws.Cells.Item(4, 35).Value = 8.42
$ws.Cells.Item(4, 36).Value = 14735240

# Row 5
$ws.Cells.Item(5, 4).Value = 9489
$ws.Cells.Item(5, 5).Value = 342
$ws.Cells.Item(5, 6).Value = 342
$ws.Cells.Item(5, 7).Value = 1772
$ws.Cells.Item(5, 8).Value = 1753
$ws.Cells.Item(5, 9).Value = 1696
$ws.Cells.Item(5, 10).Value = 57
$ws.Cells.Item(5, 11).Value = 13113
$ws.Cells.Item(5, 12).Value = 8184
$ws.Cells.Item(5, 13).Value = 4930
$ws.Cells.Item(5, 14).Value = 2875
$ws.Cells.Item(5, 15).Value = 2055
$ws.Cells.Item(5, 16).Value = 77
$ws.Cells.Item(5, 17).Value = 432
$ws.Cells.Item(5, 18).Value = -399
$ws.Cells.Item(5, 19).Value = -177
$ws.Cells.Item(5, 20).Value = 378
$ws.Cells.Item(5, 21).Value = 53
$ws.Cells.Item(5, 22).Value = 4403
$ws.Cells.Item(5, 23).Value = 3.6
$ws.Cells.Item(5, 24).Value = 18.48
$ws.Cells.Item(5, 25).Value = 58.11
$ws.Cells.Item(5, 26).Value = 14.09
$ws.Cells.Item(5, 27).Value = 166
$ws.Cells.Item(5, 28).Value = 8028.64
$ws.Cells.Item(5, 29).Value = 13844
$ws.Cells.Item(5, 30).Value = 1.09
$ws.Cells.Item(5, 31).Value = 18982
$ws.Cells.Item(5, 32).Value = 0.8
$ws.Cells.Item(5, 33).Value = 100
$ws.Cells.Item(5, 34).Value = 0.66
$ws.Cells.Item(5, 35).Value = 0.89
$ws.Cells.Item(5, 36).Value = 14847347

# Row 6
$ws.Cells.Item(6, 4).Value = 10921
$ws.Cells.Item(6, 5).Value = 415
$ws.Cells.Item(6, 6).Value = 415
$ws.Cells.Item(6, 7).Value = 173
$ws.Cells.Item(6, 8).Value = 108
$ws.Cells.Item(6, 9).Value = 47
$ws.Cells.Item(6, 11).Value = 12865
$ws.Cells.Item(6, 12).Value = 7928
$ws.Cells.Item(6, 13).Value = 4937
$ws.Cells.Item(6, 14).Value = 2895
$ws.Cells.Item(6, 16).Value = 77
$ws.Cells.Item(6, 17).Value = 811
$ws.Cells.Item(6, 18).Value = -635
$ws.Cells.Item(6, 19).Value = -19
$ws.Cells.Item(6, 20).Value = 638
$ws.Cells.Item(6, 21).Value = 173
$ws.Cells.Item(6, 22).Value = 4592
$ws.Cells.Item(6, 23).Value = 3.8
$ws.Cells.Item(6, 24).Value = 0.99
$ws.Cells.Item(6, 25).Value = 1.63
$ws.Cells.Item(6, 26).Value = 0.83
$ws.Cells.Item(6, 27).Value = 160.59
$ws.Cells.Item(6, 28).Value = 4525.92
$ws.Cells.Item(6, 29).Value = 304
$ws.Cells.Item(6, 30).Value = 37.13
$ws.Cells.Item(6, 31).Value = 19117
$ws.Cells.Item(6, 32).Value = 0.59
$ws.Cells.Item(6, 33).Value = 100
$ws.Cells.Item(6, 34).Value = 0.88
$ws.Cells.Item(6, 35).Value = 32.27
$ws.Cells.Item(6, 36).Value = 14847347

# Remove D:AI cells entirely for rows 7-9 (only A,B,C remain)
$ws.Range("D7:AI9").ClearContents()
